# Commit: "add the field dosage"
# Adds a new attribute row ("Concentración") to the Atributos sheet's
# product-attributes table, and fills in the previously-empty
# Longitud / Valor por defecto / Multivaluado / Derivado columns for
# every attribute row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Atributos": insert a new attribute row (dosage/Concentración)
# right before the "Brand" row, then fill in the new columns.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Atributos")

# Insert a new row above row 9 ("Brand"); everything below shifts down.
$ws.Range("A9").EntireRow.Insert()

# Seed the new row with the same formatting as the row above it
# (Description), then overwrite its contents below.
$ws.Range("A8:K8").Copy($ws.Range("A9:K9"))
$ws.Rows.Item(9).RowHeight = 13.8

# New row 9: the "dosage" (Concentración) attribute.
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "Concentración "
$ws.Range("D9").Value = "NO"
$ws.Range("E9").Value = "text"
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = "-"
$ws.Range("H9").Value = "NO"
$ws.Range("I9").Value = "NO"
$ws.Range("J9").Value = "NO"
$ws.Range("K9").Value = ""

# Row 6 - Id
$ws.Range("F6").Value = "integerMaxValue"
$ws.Range("G6").Value = "-"
$ws.Range("I6").Value = "NO"
$ws.Range("J6").Value = "NO"

# Row 7 - Name
$ws.Range("F7").Value = 100
$ws.Range("G7").Value = "-"
$ws.Range("I7").Value = "NO"
$ws.Range("J7").Value = "NO"

# Row 8 - Description
$ws.Range("F8").Value = 255
$ws.Range("G8").Value = "-"
$ws.Range("I8").Value = "NO"
$ws.Range("J8").Value = "NO"

# Row 10 - Brand (was row 9)
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = "-"
$ws.Range("I10").Value = "NO"
$ws.Range("J10").Value = "NO"

# Row 11 - Price (was row 10)
$ws.Range("F11").Value = "real max Value"
$ws.Range("G11").Value = "-"
$ws.Range("I11").Value = "NO"
$ws.Range("J11").Value = "NO"

# Row 12 - Stock (was row 11)
$ws.Range("F12").Value = "integerMaxValue"
$ws.Range("G12").Value = "-"
$ws.Range("I12").Value = "NO"
$ws.Range("J12").Value = "NO"

# Row 13 - Image (was row 12)
$ws.Range("F13").Value = 255
$ws.Range("G13").Value = "-"
$ws.Range("I13").Value = "NO"
$ws.Range("J13").Value = "NO"

# ---------------------------------------------------------------
# Sheet "Relaciones": a blank row 14 (default style) now exists
# between the table body and the footer row.
# ---------------------------------------------------------------
$wsRel = $wb.Worksheets.Item("Relaciones")
$wsRel.Range("A2:E2").Copy($wsRel.Range("A14:E14"))
